$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First mission-statement paragraph: the run containing "...judgement,"
#    is split across a (now-obsolete) "_GoBack" bookmark from the run
#    containing " self-reliance, independence, ambition and compassion.".
#    Re-typing the whole sentence as one contiguous Find/Replace merges the
#    two runs back into a single run and drops the bookmark in the process.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The core values Wribbenhall School seeks to promote in its children and adults are positive self-esteem; confidence in their own judgement, self-reliance, independence, ambition and compassion.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The core values Wribbenhall School seeks to promote in its children and adults are positive self-esteem; confidence in their own judgement, self-reliance, independence, ambition and compassion.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Second paragraph: "SEMH" is wrapped in proofErr spell-check markers
#    splitting the sentence into three runs. Re-typing it whole merges the
#    runs into one and clears the proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "To support our children with SEMH, anxiety and school refusal, we aim to nurture and inspire, happy re-engaged children, who are confident to engage with life to their maximum potential.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To support our children with SEMH, anxiety and school refusal, we aim to nurture and inspire, happy re-engaged children, who are confident to engage with life to their maximum potential.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Third paragraph: "curriculum" becomes "curricula", the word "a" before
#    "highly bespoke" is dropped, and the sentence is re-split into six runs.
#    Replace the whole paragraph's contents with OOXML that already has the
#    desired run boundaries.
# ---------------------------------------------------------------------------
$thirdParaText = "We seek to promote this through a highly bespoke, pupil-led curriculum encompassing"
$found3 = $d.Content.Find.Execute($thirdParaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not locate third mission-statement paragraph"
}
$p3 = $d.Paragraphs.Item(10).Range
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t>We seek to promote this through</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve">highly bespoke, pupil-led </w:t></w:r><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t>curricula</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve"> encompassing both academic, therapeutic and nurture support that embeds within the school’s core values.”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) Footer: "25 September 2019" becomes "15" + superscript "th" + " April 2021"
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("25 September 2019", $true, $false, $false, $false, $false, $true, 1, $false, "15th April 2021", 2) | Out-Null

$thRange = $footer.Range
$foundTh = $thRange.Find.Execute("th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundTh) {
    $thRange.Font.Superscript = $true
}
